$d = $word.ActiveDocument

# Locate paragraphs by their text content, since indices shift as we edit.
function Get-ParaByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# 1. Delete the entire "Folículo ovárico" paragraph (including its paragraph mark).
$folPara = Get-ParaByText "Folículo ovárico"
if ($folPara -ne $null) {
    $folPara.Range.Delete()
}

# 2. Move the "_GoBack" bookmark from the trailing empty paragraph to just before
#    the "Órganos reproductores" run (Bookmarks.Add relocates a same-named bookmark).
$orgPara = Get-ParaByText "Órganos reproductores"
$bmRange = $orgPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Remove the "Implantación" text from its paragraph, leaving the paragraph empty.
$implantPara = Get-ParaByText "Implantación"
$implantRange = $implantPara.Range.Duplicate
[void]$implantRange.MoveEnd(1, -1)
$implantRange.Delete()
